# Extend the annual SUM formula in C12 to include row 11 (previously C4:C9, now C4:C11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Formula = "=SUM(C4:C11)"

# Update the saved selection to the cell that was last active when the file was saved
$ws.Range("C12").Select()
